$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 23921
$ws.Range("I94").Value = 23921
$ws.Range("K94").Value = 23921
$ws.Range("M94").Value = -23470

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8077.3687
$ws.Range("I2").Value = 12770
$ws.Range("J2").Value = 2863.3333
$ws.Range("K2").Value = 12770
$ws.Range("L2").Value = 2863.3333
$ws.Range("M2").Value = -12657
$ws.Range("N2").Value = -3089.3333
$ws.Range("H10").Value = 450025000
$ws.Range("I10").Value = 450025000
$ws.Range("K10").Value = 450025000
$ws.Range("M10").Value = -450024830
$ws.Range("H116").Value = 8077.3687
$ws.Range("I116").Value = 12770
$ws.Range("J116").Value = 2863.3333
$ws.Range("K116").Value = 12770
$ws.Range("L116").Value = 2863.3333
$ws.Range("M116").Value = -10476
$ws.Range("N116").Value = -7451.3333
$ws.Range("H122").Value = 2061.4167
$ws.Range("I122").Value = 2081.889
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6245.667
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3795.667
$ws.Range("N122").Value = -10900

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8077.3687
$ws.Range("I3").Value = 12770
$ws.Range("J3").Value = 2863.3333
$ws.Range("K3").Value = 12770
$ws.Range("L3").Value = 2863.3333
$ws.Range("M3").Value = -12656
$ws.Range("N3").Value = -3091.3333

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000000000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1000000000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1000000000
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -1000000224
$ws.Range("H7").Value = 49.714287
$ws.Range("I7").Value = 29.75
$ws.Range("K7").Value = 29.75
$ws.Range("M7").Value = 83.25
$ws.Range("H22").Value = 137.2
$ws.Range("I22").Value = 136.76923
$ws.Range("J22").Value = 140
$ws.Range("K22").Value = 136.76923
$ws.Range("L22").Value = 140
$ws.Range("M22").Value = 213.23077
$ws.Range("N22").Value = -840
$ws.Range("H99").Value = 4465185
$ws.Range("I99").Value = 5682632.5
$ws.Range("J99").Value = 1209.3334
$ws.Range("K99").Value = 5682632.5
$ws.Range("L99").Value = 1209.3334
$ws.Range("M99").Value = -5681134.5
$ws.Range("N99").Value = -4205.3334
$ws.Range("H122").Value = 2474.818
$ws.Range("I122").Value = 1603.1428
$ws.Range("J122").Value = 4000.25
$ws.Range("K122").Value = 4809.428400000001
$ws.Range("L122").Value = 12000.75
$ws.Range("M122").Value = -2359.428400000001
$ws.Range("N122").Value = -16900.75
$ws.Range("H126").Value = 4465185
$ws.Range("I126").Value = 5682632.5
$ws.Range("J126").Value = 1209.3334
$ws.Range("K126").Value = 17047897.5
$ws.Range("L126").Value = 3628.0002
$ws.Range("M126").Value = -17045427.5
$ws.Range("N126").Value = -8568.0002

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1450
$ws.Range("I51").Value = 1450
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 4350
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -3890
$ws.Range("N51").Value = ""
$ws.Range("H68").Value = 12610.3
$ws.Range("I68").Value = 30376
$ws.Range("J68").Value = 766.5
$ws.Range("K68").Value = 91128
$ws.Range("L68").Value = 2299.5
$ws.Range("M68").Value = -90317
$ws.Range("N68").Value = -3921.5
$ws.Range("H71").Value = 12610.3
$ws.Range("I71").Value = 30376
$ws.Range("J71").Value = 766.5
$ws.Range("K71").Value = 273384
$ws.Range("L71").Value = 6898.5
$ws.Range("M71").Value = -269328
$ws.Range("N71").Value = -15010.5
$ws.Range("H88").Value = 5250
$ws.Range("J88").Value = 5250
$ws.Range("L88").Value = 15750
$ws.Range("N88").Value = -16606
$ws.Range("H91").Value = 5250
$ws.Range("J91").Value = 5250
$ws.Range("L91").Value = 15750
$ws.Range("N91").Value = -18714
$ws.Range("H113").Value = 17858048
$ws.Range("J113").Value = 18519436
$ws.Range("L113").Value = 55558308
$ws.Range("N113").Value = -55562648
$ws.Range("H121").Value = 751.6667
$ws.Range("I121").Value = 262.5
$ws.Range("J121").Value = 996.25
$ws.Range("K121").Value = 787.5
$ws.Range("L121").Value = 2988.75
$ws.Range("M121").Value = 522.5
$ws.Range("N121").Value = -5608.75
$ws.Range("H129").Value = 2126.35
$ws.Range("I129").Value = 2441.4285
$ws.Range("J129").Value = 1956.6923
$ws.Range("K129").Value = 7324.2855
$ws.Range("L129").Value = 5870.0769
$ws.Range("M129").Value = -2324.2855
$ws.Range("N129").Value = -15870.0769
$ws.Range("H131").Value = 1500.4407
$ws.Range("I131").Value = 486.18182
$ws.Range("J131").Value = 1732.875
$ws.Range("K131").Value = 1458.54546
$ws.Range("L131").Value = 5198.625
$ws.Range("M131").Value = 3581.45454
$ws.Range("N131").Value = -15278.625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1112425.1
$ws.Range("I122").Value = 1390006.4
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 4170019.2
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -4167569.2
$ws.Range("N122").Value = -11200
$ws.Range("H126").Value = 3047.3125
$ws.Range("I126").Value = 2550.5
$ws.Range("J126").Value = 3212.9167
$ws.Range("K126").Value = 7651.5
$ws.Range("L126").Value = 9638.750100000001
$ws.Range("M126").Value = -5181.5
$ws.Range("N126").Value = -14578.7501

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3329.4783
$ws.Range("I122").Value = 2197.1667
$ws.Range("J122").Value = 3729.1177
$ws.Range("K122").Value = 6591.500100000001
$ws.Range("L122").Value = 11187.3531
$ws.Range("M122").Value = -4141.500100000001
$ws.Range("N122").Value = -16087.3531

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 15100
$ws.Range("J19").Value = 3466.6667
$ws.Range("L19").Value = 3466.6667
$ws.Range("N19").Value = -3814.6667
$ws.Range("H62").Value = 16908.666
$ws.Range("I62").Value = 31072.5
$ws.Range("J62").Value = 5577.6
$ws.Range("K62").Value = 31072.5
$ws.Range("L62").Value = 5577.6
$ws.Range("M62").Value = -30448.5
$ws.Range("N62").Value = -6825.6
$ws.Range("H65").Value = 16908.666
$ws.Range("I65").Value = 31072.5
$ws.Range("J65").Value = 5577.6
$ws.Range("K65").Value = 155362.5
$ws.Range("L65").Value = 27888
$ws.Range("M65").Value = -152242.5
$ws.Range("N65").Value = -34128
$ws.Range("H126").Value = 48961.906
$ws.Range("I126").Value = 48961.906
$ws.Range("K126").Value = 146885.718
$ws.Range("M126").Value = -151684.4
